$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently stores the text "R40" (a shared string). The edit
# replaces that text with the text "1" -- still a text value (not the
# number 1), while leaving B11's existing cell style (s="23") untouched.
#
# A plain `$ws.Range("B11").Value = "1"` would let Excel's usual "looks
# like a number" auto-detection turn it into the number 1, which also
# changes the cell's type from text to numeric. To force a literal text
# value without disturbing B11's own number format, stage the text in a
# scratch cell that IS formatted as Text, copy it, and paste only the
# value into B11 (PasteSpecial xlPasteValues leaves B11's own formatting
# alone). The scratch cell is then cleared so it leaves no trace.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
$scratch.Clear()
